# This script fixes the misaligned "~UC_Sets: ..." header block pairs found
# throughout the workbook. Each affected block consists of two consecutive
# rows:
#   Row N   : a single cell in column B holding "~UC_Sets: T_xx: ..."
#             -> should actually be in column A, with B left blank.
#   Row N+1 : column A holds a marker such as "~UC_T" / "~TFM_INS" / "~TFM_UPD"
#             and column B holds "~UC_Sets: R_x: AllRegions"
#             -> the two values were swapped: A should hold the
#             "~UC_Sets: R_x: AllRegions" text and B should hold the marker.
#
# The fix is applied to every occurrence of this pattern across the
# workbook's sheets.

$wb = $excel.ActiveWorkbook

function Fix-UcSetsBlock {
    param($ws, $rowN, $rowN1)

    # --- Row N: move the value from column B to column A ---
    $bVal = $ws.Cells.Item($rowN, 2).Value2
    $ws.Cells.Item($rowN, 1).Value = $bVal
    $ws.Cells.Item($rowN, 2).Value = $null

    # --- Row N+1: swap the values held in columns A and B ---
    $aVal = $ws.Cells.Item($rowN1, 1).Value2
    $bVal2 = $ws.Cells.Item($rowN1, 2).Value2
    $ws.Cells.Item($rowN1, 1).Value = $bVal2
    $ws.Cells.Item($rowN1, 2).Value = $aVal
}

# Sheet name, row N, row N+1 for every affected block, in document order.
$blocks = @(
    @{ Sheet = "Cars";           RowN = 1;  RowN1 = 2  },
    @{ Sheet = "Cars";           RowN = 7;  RowN1 = 8  },
    @{ Sheet = "Cars_2020";      RowN = 1;  RowN1 = 2  },
    @{ Sheet = "CCS+h2";         RowN = 1;  RowN1 = 2  },
    @{ Sheet = "CH_RH";          RowN = 1;  RowN1 = 2  },
    @{ Sheet = "IND_fuels";      RowN = 1;  RowN1 = 2  },
    @{ Sheet = "Power_sector";   RowN = 1;  RowN1 = 2  },
    @{ Sheet = "Power_sector";   RowN = 10; RowN1 = 11 },
    @{ Sheet = "Power_sector";   RowN = 17; RowN1 = 18 },
    @{ Sheet = "Thermal_gencap"; RowN = 1;  RowN1 = 2  },
    @{ Sheet = "Thermal_gencap"; RowN = 7;  RowN1 = 8  },
    @{ Sheet = "TRA_Policy";     RowN = 1;  RowN1 = 2  }
)

foreach ($block in $blocks) {
    $ws = $wb.Worksheets.Item($block.Sheet)
    Fix-UcSetsBlock $ws $block.RowN $block.RowN1
}
